$d = $word.ActiveDocument

$replacements = @(
    @{ Old = "PARTNER & SENIOR DATA ARCHITECT - Siege Analytics, Washington, DC | January 2014 – Present";
       New = "PARTNER - Siege Analytics, Washington, DC | January 2014 – Present" },
    @{ Old = "PRINCIPAL TECHNICAL ARCHITECT - Clarity and Rigour, Washington, DC | 2012 – 2014";
       New = "DATA PRODUCTS MANAGER - Helm/Murmuration, Washington, DC | 2012 – 2014" },
    @{ Old = "DIRECTOR OF TECHNOLOGY - Helm, Washington, DC | 2010 – 2012";
       New = "SOFTWARE ENGINEER - Mautinoa Technologies, Washington, DC | 2010 – 2012" },
    @{ Old = "SENIOR TECHNICAL ANALYST - GSD&M, Austin, TX | 2008 – 2010";
       New = "SENIOR ANALYST - Myers Research, Washington, DC | 2008 – 2010" },
    @{ Old = "TECHNICAL COORDINATOR - Progressive Change Campaign Committee, Washington, DC | 2006 – 2008";
       New = "RESEARCH DIRECTOR - Progressive Change Campaign Committee, Washington, DC | 2006 – 2008" },
    @{ Old = "TECHNOLOGY MANAGER - The Praxis Project, Oakland, CA | 2002 – 2004";
       New = "INTERIM TECHNOLOGY MANAGER - The Praxis Project, Oakland, CA | 2002 – 2004" },
    @{ Old = "TECHNICAL COORDINATOR - The Feldman Group, Washington, DC | 2000 – 2001";
       New = "FIELD DIRECTOR - The Feldman Group, Washington, DC | 2000 – 2001" }
)

foreach ($r in $replacements) {
    $range = $d.Content
    $found = $range.Find.Execute($r.Old, $true, $false, $false, $false, $false, $true, 1, $false, $r.New, 2)
    if (-not $found) {
        Write-Host "NOT FOUND: $($r.Old)"
    } else {
        Write-Host "Replaced: $($r.Old) -> $($r.New)"
    }
}

$d.Save()
